# Apply Kujata_Profits.xlsx market-price/profit recompute (scheduled runner update).
# Each leve row has raw (non-formula) cached values for currentAveragePrice(NQ/HQ)
# and LevePrice/LeveProfit columns (H-N); we overwrite them with the refreshed values
# and drop any profit cell that no longer applies (HQ/NQ branch collapsed to 0).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1162.5
$ws.Range("J17").Value = 1285.7142
$ws.Range("L17").Value = 3857.1426
$ws.Range("N17").Value = -4193.142599999999

$ws.Range("H28").Value = 1484.8148
$ws.Range("I28").Value = 1764.2273
$ws.Range("K28").Value = 1764.2273
$ws.Range("M28").Value = -1279.2273

$ws.Range("H33").Value = 405.125
$ws.Range("I33").Value = 379.65384
$ws.Range("J33").Value = 515.5
$ws.Range("K33").Value = 379.65384
$ws.Range("L33").Value = 515.5
$ws.Range("M33").Value = -150.65384
$ws.Range("N33").Value = -973.5

$ws.Range("H54").Value = 2200
$ws.Range("I54").Value = 2200
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 2200
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -1714
$ws.Range("N54").ClearContents()

$ws.Range("H64").Value = 3501.7
$ws.Range("I64").Value = 3498.2856
$ws.Range("J64").Value = 3509.6667
$ws.Range("K64").Value = 3498.2856
$ws.Range("L64").Value = 3509.6667
$ws.Range("M64").Value = -3250.2856
$ws.Range("N64").Value = -4005.6667

$ws.Range("H67").Value = 3501.7
$ws.Range("I67").Value = 3498.2856
$ws.Range("J67").Value = 3509.6667
$ws.Range("K67").Value = 3498.2856
$ws.Range("L67").Value = 3509.6667
$ws.Range("M67").Value = -2640.2856
$ws.Range("N67").Value = -5225.6667

$ws.Range("H80").Value = 7014.5713
$ws.Range("I80").Value = 22501
$ws.Range("J80").Value = 820
$ws.Range("K80").Value = 67503
$ws.Range("L80").Value = 2460
$ws.Range("M80").Value = -66505
$ws.Range("N80").Value = -4456

$ws.Range("H83").Value = 7014.5713
$ws.Range("I83").Value = 22501
$ws.Range("J83").Value = 820
$ws.Range("K83").Value = 202509
$ws.Range("L83").Value = 7380
$ws.Range("M83").Value = -197517
$ws.Range("N83").Value = -17364

$ws.Range("H98").Value = 2748.3333
$ws.Range("I98").Value = 3033
$ws.Range("J98").Value = 1325
$ws.Range("K98").Value = 3033
$ws.Range("L98").Value = 1325
$ws.Range("M98").Value = -1535
$ws.Range("N98").Value = -4321

$ws.Range("H122").Value = 2748.3333
$ws.Range("I122").Value = 3033
$ws.Range("J122").Value = 1325
$ws.Range("K122").Value = 9099
$ws.Range("L122").Value = 3975
$ws.Range("M122").Value = -6649
$ws.Range("N122").Value = -8875

$ws.Range("H126").Value = 40000
$ws.Range("J126").Value = 40000
$ws.Range("L126").Value = 40000
$ws.Range("N126").Value = -49880

$ws.Range("H129").Value = 809.7619
$ws.Range("I129").Value = 503.22223
$ws.Range("J129").Value = 860.85187
$ws.Range("K129").Value = 1509.66669
$ws.Range("L129").Value = 2582.55561
$ws.Range("M129").Value = 3490.33331
$ws.Range("N129").Value = -12582.55561

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 794.2105
$ws.Range("I2").Value = 617
$ws.Range("J2").Value = 991.1111
$ws.Range("K2").Value = 617
$ws.Range("L2").Value = 991.1111
$ws.Range("M2").Value = -504
$ws.Range("N2").Value = -1217.1111

$ws.Range("H32").Value = 5331.8794
$ws.Range("I32").Value = 4287.4526
$ws.Range("K32").Value = 4287.4526
$ws.Range("M32").Value = -4000.4526

$ws.Range("H45").Value = 1024.069
$ws.Range("I45").Value = 961.46155
$ws.Range("J45").Value = 1566.6666
$ws.Range("K45").Value = 961.46155
$ws.Range("L45").Value = 1566.6666
$ws.Range("M45").Value = -584.46155
$ws.Range("N45").Value = -2320.6666

$ws.Range("H116").Value = 794.2105
$ws.Range("I116").Value = 617
$ws.Range("J116").Value = 991.1111
$ws.Range("K116").Value = 617
$ws.Range("L116").Value = 991.1111
$ws.Range("M116").Value = 1677
$ws.Range("N116").Value = -5579.1111

$ws.Range("H122").Value = 1832
$ws.Range("I122").Value = 1832
$ws.Range("K122").Value = 5496
$ws.Range("M122").Value = -3046

$ws.Range("H132").Value = 3347.6
$ws.Range("I132").Value = 2856
$ws.Range("K132").Value = 8568
$ws.Range("M132").Value = -6038

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 794.2105
$ws.Range("I3").Value = 617
$ws.Range("J3").Value = 991.1111
$ws.Range("K3").Value = 617
$ws.Range("L3").Value = 991.1111
$ws.Range("M3").Value = -503
$ws.Range("N3").Value = -1219.1111

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H107").Value = 1179.8
$ws.Range("I107").Value = 880.9167
$ws.Range("J107").Value = 2375.3333
$ws.Range("K107").Value = 880.9167
$ws.Range("L107").Value = 2375.3333
$ws.Range("M107").Value = 1039.0833
$ws.Range("N107").Value = -6215.3333

$ws.Range("H134").Value = 6138.773
$ws.Range("I134").Value = 1131.6666
$ws.Range("J134").Value = 16868.285
$ws.Range("K134").Value = 3394.9998
$ws.Range("L134").Value = 50604.855
$ws.Range("M134").Value = -859.9998000000001
$ws.Range("N134").Value = -55674.855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 503.42856
$ws.Range("I105").Value = 454
$ws.Range("K105").Value = 454
$ws.Range("M105").Value = 1293

$ws.Range("H132").Value = 1881.1666
$ws.Range("I132").Value = 1545.875
$ws.Range("K132").Value = 4637.625
$ws.Range("M132").Value = -2107.625

$ws.Range("H134").Value = 1684.1923
$ws.Range("J134").Value = 2338
$ws.Range("L134").Value = 7014
$ws.Range("N134").Value = -12084

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 541.5714
$ws.Range("I5").Value = 517.75
$ws.Range("J5").Value = 573.3333
$ws.Range("K5").Value = 1553.25
$ws.Range("L5").Value = 1719.9999
$ws.Range("M5").Value = -1441.25
$ws.Range("N5").Value = -1943.9999

$ws.Range("H20").Value = 237.5
$ws.Range("I20").Value = 200
$ws.Range("K20").Value = 600
$ws.Range("M20").Value = -373

$ws.Range("H94").Value = 3940.5
$ws.Range("J94").Value = 4050
$ws.Range("L94").Value = 12150
$ws.Range("N94").Value = -13502

$ws.Range("H107").Value = 11410.3
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 12567
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 37701
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -41541

$ws.Range("H131").Value = 12348204
$ws.Range("I131").Value = 142857440
$ws.Range("J131").Value = 2735.7026
$ws.Range("K131").Value = 428572320
$ws.Range("L131").Value = 8207.1078
$ws.Range("M131").Value = -428567280
$ws.Range("N131").Value = -18287.1078

$ws.Range("H132").Value = 1501
$ws.Range("I132").Value = 1001
$ws.Range("J132").Value = 2001
$ws.Range("K132").Value = 9009
$ws.Range("L132").Value = 18009
$ws.Range("M132").Value = -6479
$ws.Range("N132").Value = -23069

$ws.Range("H135").Value = 541.5714
$ws.Range("I135").Value = 517.75
$ws.Range("J135").Value = 573.3333
$ws.Range("K135").Value = 4659.75
$ws.Range("L135").Value = 5159.9997
$ws.Range("M135").Value = -2124.75
$ws.Range("N135").Value = -10229.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 510.2381
$ws.Range("I102").Value = 496.05264
$ws.Range("J102").Value = 645
$ws.Range("K102").Value = 496.05264
$ws.Range("L102").Value = 645
$ws.Range("M102").Value = 1125.94736
$ws.Range("N102").Value = -3889

$ws.Range("H122").Value = 1864.56
$ws.Range("I122").Value = 1897.8125
$ws.Range("J122").Value = 1805.4445
$ws.Range("K122").Value = 5693.4375
$ws.Range("L122").Value = 5416.333500000001
$ws.Range("M122").Value = -3243.4375
$ws.Range("N122").Value = -10316.3335

$ws.Range("H126").Value = 1963.5652
$ws.Range("I126").Value = 1707.2941
$ws.Range("J126").Value = 2689.6667
$ws.Range("K126").Value = 5121.8823
$ws.Range("L126").Value = 8069.000100000001
$ws.Range("M126").Value = -2651.8823
$ws.Range("N126").Value = -13009.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3739.9285
$ws.Range("I46").Value = 482.7143
$ws.Range("K46").Value = 482.7143
$ws.Range("M46").Value = -294.7143

$ws.Range("H55").Value = 526.85
$ws.Range("I55").Value = 616.61536
$ws.Range("J55").Value = 360.14285
$ws.Range("K55").Value = 616.61536
$ws.Range("L55").Value = 360.14285
$ws.Range("M55").Value = -443.61536
$ws.Range("N55").Value = -706.14285

$ws.Range("H93").Value = 569.9
$ws.Range("I93").Value = 564.1429000000001
$ws.Range("J93").Value = 583.3333
$ws.Range("K93").Value = 564.1429000000001
$ws.Range("L93").Value = 583.3333
$ws.Range("M93").Value = 683.8570999999999
$ws.Range("N93").Value = -3079.3333

$ws.Range("H132").Value = 3058.8
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 3176.4443
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 9529.332900000001
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -14589.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 22728848
$ws.Range("I122").Value = 22728848
$ws.Range("K122").Value = 68186544
$ws.Range("M122").Value = -68184094

$ws.Range("H133").Value = 37871.668
$ws.Range("J133").Value = 37871.668
$ws.Range("L133").Value = 37871.668
$ws.Range("N133").Value = -47991.668
